$wb = $excel.ActiveWorkbook

# --- Sheet: Linear ---
# Update mu, B, sig2 parameter estimates, and the abs_epsi_autocorr series (B5)
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -0.0004812284099154532
$wsLinear.Range("B3").Value = -0.003608829674642622
$wsLinear.Range("B4").Value = 0.02966370961644622
$wsLinear.Range("B5").Value = "[1.0, 0.2151549308327078, 0.07816655677432746, 0.07930948181075301, 0.046878059965325874, 0.06235897953964716, 0.21819458796998653, 0.37235869009907596, 0.2149235259311273, 0.06386939784605382, 0.024616406753886853, 0.04026214170303388, 0.06893680552945473, 0.20973289557518948, 0.3632339344392674, 0.22332091506090718, 0.027508033972460287, 0.039149997087396964, 0.035311814440135786, 0.049619895437818175]"

# --- Sheet: NonLinear ---
# Update sig2_0, mu_1, B_1, sig2_1, c, p parameter estimates, and the abs_epsi_autocorr series (B10)
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = 0.003005388459806029
$wsNonLinear.Range("B5").Value = 0.002984055386783263
$wsNonLinear.Range("B6").Value = 0.02552841863866622
$wsNonLinear.Range("B7").Value = -0.001497563604437836
$wsNonLinear.Range("B8").Value = -0.003068042538672286
$wsNonLinear.Range("B9").Value = 0.03406559366896734
$wsNonLinear.Range("B10").Value = "[1.0, 0.21501292045492706, 0.07917993512678827, 0.08034794535101615, 0.04789143869140638, 0.06317667931439425, 0.21826211596418207, 0.37162834915948617, 0.2147414345914138, 0.06481067548774962, 0.025314991838125758, 0.04097815421347617, 0.06961050446068438, 0.20952404817173312, 0.36244279649568034, 0.22347450761702337, 0.02828018752108814, 0.039840797173624896, 0.0361681579495062, 0.05020786007456066]"
